$wb = $excel.ActiveWorkbook

# Update the Initial Time year value on the "IT" sheet
$itSheet = $wb.Worksheets.Item("IT")
$itSheet.Range("B2").Value = 2020

# Make "IT" the active/selected sheet (was "About") and set its selection to B3
$itSheet.Activate()
$itSheet.Range("B3").Select()
